$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column date serial from 45208 to 45212 for rows 2-10
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value2 = 45212
    }
}
